$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original values referenced "ERR_NOT_FRIENDS" for commands that depend on
# the trust relationship (download, msg). Since the underlying file that backs
# that relationship is actually the "trusted" list (not a generic "friends"
# concept), the error code text is corrected to "ERR_NOT_TRUSTED" in both the
# detailed table (rows 7-27) and the summary table (rows 30-50).

$ws.Range("E21").Value = "ERR_YOURSELF or ERR_NOT_FOUND or ERR_NOT_TRUSTED or ERR_NOT_REGISTERED"
$ws.Range("E23").Value = "OP_ERROR or ERR_NOT_TRUSTED or ERR_NOT_REGISTERED"
$ws.Range("E44").Value = "ERR_YOURSELF or ERR_NOT_FOUND or ERR_NOT_TRUSTED or ERR_NOT_REGISTERED"
$ws.Range("E46").Value = "OP_ERROR or ERR_NOT_TRUSTED or ERR_NOT_REGISTERED"

# Reflect the saved view state: the last active selection moved to the
# second (summary) table, landing on G46.
$ws.Range("G46").Select() | Out-Null
